# Insert a new data row at row 216 (pushing the existing rows 216..308 down
# to 217..309) and populate the new row with the weekly price-report entry
# for Acelga / Femacal de La Calera (Coquimbo) that was added upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 216..308 down by one row, leaving a blank row 216.
$ws.Rows.Item(216).Insert()

# Fill in the new row 216 with the new record's values.
$ws.Cells.Item(216, 1).Value  = 3
$ws.Cells.Item(216, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(216, 3).Value  = "Coquimbo"
$ws.Cells.Item(216, 4).Value  = 44636
$ws.Cells.Item(216, 5).Value  = 5
$ws.Cells.Item(216, 6).Value  = 100112009
$ws.Cells.Item(216, 7).Value  = "Acelga"
$ws.Cells.Item(216, 8).Value  = "Sin especificar"
$ws.Cells.Item(216, 9).Value  = "Primera"
$ws.Cells.Item(216, 10).Value = 195
$ws.Cells.Item(216, 11).Value = 4000
$ws.Cells.Item(216, 12).Value = 4500
$ws.Cells.Item(216, 13).Value = 4282
$ws.Cells.Item(216, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(216, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(216, 16).Value = 714
$ws.Cells.Item(216, 17).Value = 6
$ws.Cells.Item(216, 18).Value = "Hortaliza"
